# 卫斯理仓位统计 - add the 20191103 snapshot sheet.
#
# Strategy: duplicate the most recent existing sheet (20191029) so the new
# sheet inherits its column widths, number formats and the D/E "shared"
# formula block, rename it to 20191103, then overwrite the data cells with
# the new day's numbers. Excel/the engine recalculates E1:E13, F14, B15,
# C15 and B24/C24 automatically because they're formulas.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("20191029")
$template.Copy($null, $template) | Out-Null
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "20191103"

# ---- position-bucket table (rows 1-14) ----
$ws.Range("B1").Value = 812
$ws.Range("B2").Value = 336
$ws.Range("B3").Value = 388
$ws.Range("B4").Value = 449
$ws.Range("B5").Value = 487
$ws.Range("C5").Value = 0.03
$ws.Range("B6").Value = 558
$ws.Range("B7").Value = 779
$ws.Range("B8").Value = 650
$ws.Range("B9").Value = 753
$ws.Range("B10").Value = 999
$ws.Range("C10").Value = 0.08
$ws.Range("B11").Value = 1463
$ws.Range("C11").Value = 0.11
$ws.Range("B12").Value = 1535
$ws.Range("B13").Value = 832
$ws.Range("C13").Value = 0.06
$ws.Range("B14").Value = 2177
$ws.Range("C14").Value = 0.17

# ---- sentiment poll (rows 20-23) ----
# previously "看多" / "看空 (已选)" / "看平" / call-out text; this round the
# "看多" choice is the selected one.
$ws.Range("A20").Value = "看多 (已选)"
$ws.Range("B20").Value = 5942
$ws.Range("C20").Value = 0.48

$ws.Range("A21").Value = "看空"
$ws.Range("B21").Value = 2322
$ws.Range("C21").Value = 0.19

$ws.Range("B22").Value = 1728
$ws.Range("C22").Value = 0.14

$ws.Range("A23").Value = "我是来给卫斯理打Call的~"
$ws.Range("B23").Value = 2225
$ws.Range("C23").Value = 0.18

# ---- view bookkeeping ----
# New sheet becomes the active tab, selection moves to D14; old sheet keeps
# its previous selection (M17) but is no longer the selected tab.
$ws.Range("D14").Select() | Out-Null

$old = $wb.Worksheets.Item("20191029")
$old.Range("M17").Select() | Out-Null
$ws.Activate()

# Scroll both tabs to match the authored viewport (best-effort; some
# headless hosts don't persist viewport scroll position).
$win = $excel.ActiveWindow
$old.Activate()
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Activate()
$win.ScrollRow = 6
$win.ScrollColumn = 1
$win.ScrollWorkbookTabs(2) | Out-Null
